$d = $word.ActiveDocument

function Replace-Text($old, $new, $wholeWord) {
    $result = $d.Content.Find.Execute($old, $true, $wholeWord, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output "WARNING: replace failed for: $old"
    }
}

# 1. Normalize font for every paragraph: TimesNewToman -> Times New Roman
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Font.Name = "Times New Roman"
}

# 2. Title
Replace-Text "Virtual Reality: An Immersive Paradigm Shift" "Understanding the Dynamic Equilibrium of Ecosystems: A Balanced Orchestra of Life" $false

# 3. Author name
Replace-Text "Varun Wadhwani" "Dr. Alexander Westwood" $false

# 4. Email address pieces
Replace-Text "varun" "westwood" $true
Replace-Text "wadhwani@emailworld" "a@eduinstitute" $false
Replace-Text "com" "org" $true

# 5. Main body paragraph - sentence by sentence, left to right
Replace-Text "Virtual Reality(VR), a captivating technology, transports users into a simulated environment, empowering them to interact with digital creations" "The intricate harmony of ecosystems lies in the delicate balance maintained between organisms and their environment" $false

Replace-Text " Its transformative properties have sparked a paradigm shift across various disciplines, from gaming and entertainment to education, healthcare, and workplace simulations" " Life thrives within this dynamic equilibrium, where innumerable interactions weave a complex tapestry of interdependence" $false

Replace-Text " As VR's applications continue to expand, let's delve into its profound impact and explore the boundless possibilities it holds for shaping the future." " Every organism plays a unique role in maintaining this delicate balance, contributing to the overall stability and resilience of the ecosystem. To comprehend the intricate dance of life, we must delve into the fundamental principles governing these interactions." $false

Replace-Text "VR's immersive nature has revolutionized gaming and entertainment, providing an unparalleled level of engagement and escapism" "Understanding the roles of individual organisms within an ecosystem is crucial" $false

Replace-Text " Players can now step into virtual worlds, embodying characters and experiencing adventures like never before" " Each species occupies a specific ecological niche, playing a distinct role in energy flow and nutrient cycling" $false

Replace-Text " This immersive experience has also found its way into other fields, such as education and training, where simulations can provide realistic and engaging scenarios for students and professionals alike." " This interconnectedness forms a web of relationships that shape the dynamics of the ecosystem. Changes in one species can ripple through the entire system, triggering a cascade of ecological responses. By studying these interactions, scientists can unravel the intricate mechanisms that maintain equilibrium." $false

Replace-Text "Beyond entertainment and education, VR is making significant strides in healthcare" "Biodiversity, the vast array of species within an ecosystem, is crucial for maintaining ecological balance" $false

Replace-Text " It offers immersive therapies for conditions like PTSD and phobias, allowing patients to confront their fears in a controlled virtual environment" " A rich diversity of species enhances the resilience of ecosystems, allowing them to adapt to environmental changes" $false

Replace-Text " Moreover, VR is instrumental in surgical training, enabling surgeons to practice complex procedures in a risk-free environment, leading to improved surgical outcomes" " By providing a variety of habitats and resources, biodiversity ensures the survival of a wide range of organisms. Preserving biodiversity is therefore essential for the long-term stability and productivity of ecosystems" $false

# 6. Summary paragraph
Replace-Text "Virtual Reality has emerged as a pivotal technology, ushering in a new era of immersive experiences across various domains" "Ecosystems are complex systems in which organisms and their environment interact, creating a dynamic balance" $false

Replace-Text " Its applications have soared, ranging from gaming and entertainment to education, healthcare, and corporate training" " The roles of individual species, their interconnectedness, and biodiversity are key factors in maintaining this balance" $false

Replace-Text " VR's ability to transport users into digital environments has enabled novel and engaging ways of learning, healing, and simulating complex scenarios" " Understanding these interactions allows us to appreciate the delicate harmony of life and the importance of preserving biodiversity" $false

Replace-Text " As technology continues to advance, VR's impact is poised to grow exponentially, reshaping industries and transforming the way we learn, heal, and experience the world around us" " By studying ecosystems, we gain insights into the intricate dance of life, unraveling the secrets of maintaining a healthy and balanced environment" $false

# 7. Append a new empty paragraph at the very end of the document
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

Write-Output "Edit complete"
